$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header-row relabeling (row 1) ---
# F/G: Frac1_Remote/Frac1_Natural -> Frac1_DarkBlue/Frac1_LightBlue
$ws.Range("F1").Value = "Frac1_DarkBlue"
$ws.Range("G1").Value = "Frac1_LightBlue"

# M/N/O: Frac2_Remote/Frac2_Natural/Frac2_Potential -> Frac2_DarkBlue/Frac2_LightBlue/Frac2_Grey
$ws.Range("M1").Value = "Frac2_DarkBlue"
$ws.Range("N1").Value = "Frac2_LightBlue"
$ws.Range("O1").Value = "Frac2_Grey"

# Q: Total3 -> Total3 - TEXT
$ws.Range("Q1").Value = "Total3 - TEXT"

# --- Hide helper columns D:E, I, L ---
$ws.Columns("D:E").ColumnWidth = -0.8333333333333334
$ws.Columns("D:E").Hidden = $true
$ws.Columns("I").ColumnWidth = -0.8333333333333334
$ws.Columns("I").Hidden = $true
$ws.Columns("L").ColumnWidth = -0.8333333333333334
$ws.Columns("L").Hidden = $true

# --- Apply custom "+#" number format to K2:K8 ---
$ws.Range("K2:K8").NumberFormat = "\+#"

# --- View: drop frozen/top-left scroll position, move selection ---
$ws.Range("F14").Select()
